$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 128. This shifts the existing rows 128-206
# down to 130-208, automatically growing the sheet dimension to A1:R208 and
# pushing the former last two rows (old 205, old 206) into the new rows
# 207 and 208.
$ws.Rows("128:129").Insert()

# Populate new row 128 with its data (most fields identical to the template
# row that is now at row 130, except for the changed fields below).
$ws.Range("A128").Value = 10
$ws.Range("B128").Value = "Vega Modelo de Temuco"
$ws.Range("C128").Value = "La Araucanía"
$ws.Range("D128").Value = 44518
$ws.Range("E128").Value = 9
$ws.Range("F128").Value = 100112017
$ws.Range("G128").Value = "Apio"
$ws.Range("H128").Value = "Americana (o)"
$ws.Range("I128").Value = "Primera"
$ws.Range("J128").Value = 120
$ws.Range("K128").Value = 8000
$ws.Range("L128").Value = 9000
$ws.Range("M128").Value = 8458
$ws.Range("N128").Value = "`$/docena de matas"
$ws.Range("O128").Value = "Provincia del Elquí"
$ws.Range("P128").Value = 1410
$ws.Range("Q128").Value = 6
$ws.Range("R128").Value = "Hortaliza"

# Populate new row 129 with its data.
$ws.Range("A129").Value = 10
$ws.Range("B129").Value = "Vega Modelo de Temuco"
$ws.Range("C129").Value = "La Araucanía"
$ws.Range("D129").Value = 44518
$ws.Range("E129").Value = 9
$ws.Range("F129").Value = 100112017
$ws.Range("G129").Value = "Apio"
$ws.Range("H129").Value = "Americana (o)"
$ws.Range("I129").Value = "Primera"
$ws.Range("J129").Value = 95
$ws.Range("K129").Value = 9000
$ws.Range("L129").Value = 10000
$ws.Range("M129").Value = 9474
$ws.Range("N129").Value = "`$/docena de matas"
$ws.Range("O129").Value = "Región Metropolitana"
$ws.Range("P129").Value = 1579
$ws.Range("Q129").Value = 6
$ws.Range("R129").Value = "Hortaliza"
